$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# The shakeflask data sheet is gaining a new "T6" timepoint row for each
# growth condition (rows 27, 34, 41). Those rows' OD/pH cells (E/F) were
# left with a stray one-off format; align them with the format already
# used by every other populated OD/pH cell in the table (copied from
# E2) and restore the normal 19.5pt row height used throughout the rest
# of the table (rows 27 & 34 had drifted to 17.25pt).
# ---------------------------------------------------------------------

$ws.Rows.Item(27).RowHeight = 19.5
$ws.Rows.Item(34).RowHeight = 19.5

$ws.Range("E2").Copy()
foreach ($addr in @("E27:F27", "E34:F34", "E41:F41")) {
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false

# Refresh the table's look now that the report is being drafted.
$tbl = $ws.ListObjects.Item("Table1")
$tbl.TableStyle = "TableStyleLight1"
